# Generate Report for Handoff
# Update the "latest handoff" timestamps for the bfb8ff40-... file row across
# the Overview summary sheet and the per-locale (zh-cn / de-de) detail sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("D6").Value = "2016-12-18 07:12:06"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E6").Value = "2016-03-18 07:12:03"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E6").Value = "2016-03-18 07:12:06"
